$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: rename Prophet Single -> Prophet Univariate, add XGBoost headers ---
$ws.Range("H1").Value = "Prophet Univariate"
$ws.Range("J1").Value = "XGBoost"
$ws.Range("J1").Interior.Color = $ws.Range("B1").Interior.Color
$ws.Range("L1").Value = "XGBoost Log"
$ws.Range("L1").Interior.Color = $ws.Range("B1").Interior.Color

# --- Header row 2: sub-headers for new columns (format copied from B2:C2, then values set) ---
$ws.Range("B2:C2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("J2").Value = "rmse_pred_train"
$ws.Range("K2").Value = "rmse_pred_test"
$ws.Range("L2").Value = "rmse_train_pred"
$ws.Range("M2").Value = "rmse_test_pred"

# --- Data rows 3-25: copy number-format/fill from H:I (s=2/s=3) down to J:M, then set values ---
$ws.Range("H3:I25").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("H3:I25").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# --- Update existing Prophet (H/I) values, and set new XGBoost (J/K) + XGBoost Log (L/M) values ---
$ws.Range("H3").Value = 41677.5278532
$ws.Range("I3").Value = 140116.800015394
$ws.Range("J3").Value = 0.00933748965307019
$ws.Range("K3").Value = 9762.80512208658
$ws.Range("L3").Value = 70.9374403
$ws.Range("M3").Value = 306.977191
$ws.Range("H4").Value = 44677.0084511375
$ws.Range("I4").Value = 138514.013328332
$ws.Range("J4").Value = 0.085867710476516
$ws.Range("K4").Value = 12883.0201770245
$ws.Range("L4").Value = 98.8727704
$ws.Range("M4").Value = 1579.04776
$ws.Range("H5").Value = 51885.9125897124
$ws.Range("I5").Value = 141628.281109137
$ws.Range("J5").Value = 0.010322337757313
$ws.Range("K5").Value = 5144.65704655594
$ws.Range("L5").Value = 37.1428172
$ws.Range("M5").Value = 216.147267
$ws.Range("H6").Value = 54509.4093597118
$ws.Range("I6").Value = 143262.915131289
$ws.Range("J6").Value = 0.0202528995702124
$ws.Range("K6").Value = 1711.82380124411
$ws.Range("L6").Value = 51.266564
$ws.Range("M6").Value = 53.9779187
$ws.Range("H7").Value = 63008.158693643
$ws.Range("I7").Value = 144278.248142748
$ws.Range("J7").Value = 0.0172190194731582
$ws.Range("K7").Value = 4950.62847082249
$ws.Range("L7").Value = 7.04351698
$ws.Range("M7").Value = 45.1423147
$ws.Range("H8").Value = 67631.0247515298
$ws.Range("I8").Value = 145339.34191183
$ws.Range("J8").Value = 0.165530431190444
$ws.Range("K8").Value = 291.208548331889
$ws.Range("L8").Value = 5.14762826
$ws.Range("M8").Value = 5.94672732
$ws.Range("H9").Value = 67423.530696455
$ws.Range("I9").Value = 142939.147832184
$ws.Range("J9").Value = 0.105345304445282
$ws.Range("K9").Value = 43841.3030020499
$ws.Range("L9").Value = 384.213311
$ws.Range("M9").Value = 1842.38208
$ws.Range("H10").Value = 63969.4609505365
$ws.Range("I10").Value = 143259.21010905
$ws.Range("J10").Value = 0.0618513568031147
$ws.Range("K10").Value = 43116.1672358049
$ws.Range("L10").Value = 243.125774
$ws.Range("M10").Value = 2364.65517
$ws.Range("H11").Value = 82887.6158489118
$ws.Range("I11").Value = 147094.015049272
$ws.Range("J11").Value = 0.230663872088021
$ws.Range("K11").Value = 18171.2455522458
$ws.Range("L11").Value = 394.512539
$ws.Range("M11").Value = 1671.9255
$ws.Range("H12").Value = 69751.1577741529
$ws.Range("I12").Value = 146168.6830511
$ws.Range("J12").Value = 0.248335757534345
$ws.Range("K12").Value = 17970.4994483455
$ws.Range("L12").Value = 201.724493
$ws.Range("M12").Value = 1497.95851
$ws.Range("H13").Value = 79832.0176555958
$ws.Range("I13").Value = 138041.525791273
$ws.Range("J13").Value = 0.640344075864911
$ws.Range("K13").Value = 41599.8176255249
$ws.Range("L13").Value = 291.430514
$ws.Range("M13").Value = 23808.5873
$ws.Range("H14").Value = 86543.4134100363
$ws.Range("I14").Value = 149682.07196088
$ws.Range("J14").Value = 0.0151097039465608
$ws.Range("K14").Value = 5405.22665708458
$ws.Range("L14").Value = 62.8973508
$ws.Range("M14").Value = 628.850087
$ws.Range("H15").Value = 95699.0049222422
$ws.Range("I15").Value = 152185.304366416
$ws.Range("J15").Value = 0.0342999734831001
$ws.Range("K15").Value = 737.862480923895
$ws.Range("L15").Value = 11.384246
$ws.Range("M15").Value = 6.37137517
$ws.Range("H16").Value = 99523.4474797133
$ws.Range("I16").Value = 153159.018047239
$ws.Range("J16").Value = 0.153891195391295
$ws.Range("K16").Value = 289.908032483239
$ws.Range("L16").Value = 11.0537353
$ws.Range("M16").Value = 19.0895966
$ws.Range("H17").Value = 101456.032713835
$ws.Range("I17").Value = 153714.824667251
$ws.Range("J17").Value = 0.0687585025165979
$ws.Range("K17").Value = 2341.85665388331
$ws.Range("L17").Value = 32.270944
$ws.Range("M17").Value = 400.345659
$ws.Range("H18").Value = 93194.1627231528
$ws.Range("I18").Value = 150527.409044843
$ws.Range("J18").Value = 0.181271530752019
$ws.Range("K18").Value = 7227.14989719584
$ws.Range("L18").Value = 160.638517
$ws.Range("M18").Value = 2018.30751
$ws.Range("H19").Value = 112532.685294463
$ws.Range("I19").Value = 156048.051866677
$ws.Range("J19").Value = 0.0563948304295638
$ws.Range("K19").Value = 3187.4049337862
$ws.Range("L19").Value = 6.99692635
$ws.Range("M19").Value = 45.7176747
$ws.Range("H20").Value = 113758.394926279
$ws.Range("I20").Value = 156988.158878473
$ws.Range("J20").Value = 0.018072463213376
$ws.Range("K20").Value = 865.330428199406
$ws.Range("L20").Value = 34.1400706
$ws.Range("M20").Value = 81.31174
$ws.Range("H21").Value = 114646.741232175
$ws.Range("I21").Value = 138682.120143512
$ws.Range("J21").Value = 0.00751583542656682
$ws.Range("K21").Value = 25687.8338526384
$ws.Range("L21").Value = 408.874409
$ws.Range("M21").Value = 18971.8066
$ws.Range("H22").Value = 124982.208242052
$ws.Range("I22").Value = 159028.962953407
$ws.Range("J22").Value = 0.0305684250071634
$ws.Range("K22").Value = 250.501845074337
$ws.Range("L22").Value = 6.48027083
$ws.Range("M22").Value = 1.46360958
$ws.Range("H23").Value = 115654.995300418
$ws.Range("I23").Value = 160007.047960949
$ws.Range("J23").Value = 0.250671581510491
$ws.Range("K23").Value = 9053.56668098974
$ws.Range("L23").Value = 136.852043
$ws.Range("M23").Value = 2.52431888
$ws.Range("H24").Value = 115815.900531128
$ws.Range("I24").Value = 159000.879273258
$ws.Range("J24").Value = 0.170292168318711
$ws.Range("K24").Value = 8485.92687432822
$ws.Range("L24").Value = 231.835717
$ws.Range("M24").Value = 873.296284
$ws.Range("H25").Value = 131954.995401176
$ws.Range("I25").Value = 147743.66442476
$ws.Range("J25").Value = 0.0575279248111517
$ws.Range("K25").Value = 202694.939341962
$ws.Range("L25").Value = 989.732078
$ws.Range("M25").Value = 3427.14627

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 14.3
$ws.Columns.Item(10).ColumnWidth = 15.15
$ws.Columns.Item(11).ColumnWidth = 13.02
$ws.Columns.Item(12).ColumnWidth = 13.5
$ws.Columns.Item(13).ColumnWidth = 14.5

# --- Sheet view: zoom + selection ---
$excel.ActiveWindow.Zoom = 125
$ws.Range("M3").Select()

Write-Host "edit applied"
